# "Generate Report for Archive"
# The localization run moved from hand-off to active translation, so every
# place that surfaced the old "Ready for handoff" status now reads
# "In Translation" (Overview!E2/F2 for zh-cn & de-de, and the Status column
# on each per-language sheet). The shorter status text made Excel re-fit the
# Status column narrower on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: per-language status columns (zh-cn = E, de-de = F)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-language detail sheets: Status column (C)
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Column widths shrink to fit the new, shorter status text
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
